$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(412402, "InternalFailureException - An unexpected error has occurred."),
    @(412403, "InvalidRequestException - The request is not valid."),
    @(412404, "ResourceNotFoundException - The specified resource does not exist."),
    @(412405, "ServiceUnavailableException - The service is temporarily unavailable."),
    @(412406, "ThrottlingException - The rate exceeds the limit."),
    @(412407, "TransferAlreadyCompletedException - You can't revert the certificate transfer because the transfer is already complete."),
    @(412408, "UnauthorizedException - You are not authorized to perform this operation."),
    @(412409, "LimitExceededException - The number of attached entities exceeds the limit."),
    @(412410, "MalformedPolicyException - The policy documentation is not valid."),
    @(412411, "ResourceAlreadyExistsException - The resource already exists."),
    @(412412, "VersionsLimitExceededException - The number of policy versions exceeds the limit."),
    @(412413, "InternalException - An unexpected error has occurred."),
    @(412414, "SqlParseException - The Rule-SQL expression cannot be parsed correctly."),
    @(412415, "CertificateStateException - The certificate operation is not allowed."),
    @(412416, "DeleteConflictException - You cannot delete the resource because it is attached to one or more resources."),
    @(412417, "VersionConflictException - The version of the thing is different than the version specified with the --version parameter."),
    @(412418, "CertificateValidationException - The certificate is invalid."),
    @(412419, "RegistrationCodeValidationException - The registration code is invalid."),
    @(412420, "CertificateConflictException - Unable to verify the CA certificate used to sign the device certificate you are attempting to register. This happens when you have registered `nmore than one CA certificate that has the same subject field and public key."),
    @(412421, "TransferConflictException - You cannot transfer the certificate because authorization policies are still attached.")
)

$row = 4
foreach ($entry in $data) {
    $code = $entry[0]
    $desc = $entry[1]
    $ws.Cells.Item($row, 1).Value = $code
    $ws.Cells.Item($row, 2).Value = $desc
    $row = $row + 1
}

# Row 22 (code 412420) has special formatting: wrapped text and a taller row
$ws.Range("B22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 30

$ws.Range("B26").Select()
